$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-converted to a number by Excel (losing trailing zeros / precision)
$textCells = @("D4", "D5", "D6", "D8", "D12", "D13", "D17", "D20", "D24", "D25", "D26", "D27", "D34", "D37", "D38", "D39", "D40", "D42", "D43", "D45", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

# Apply cell value changes per diff
$ws.Range("D2").Value = "63.487.50"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "3.057.59"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "590.22"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "154.53"
$ws.Range("E6").Value = "  +6.25%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("D9").Value = "3.075.59"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "37.62"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "3.564.29"
$ws.Range("E16").Value = "  -2.73%  "
$ws.Range("D17").Value = "7.21"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "63.413.48"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "3.069.09"
$ws.Range("D20").Value = "478.30"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  +3.90%  "
$ws.Range("D25").Value = "12.96"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "81.18"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "10.01"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").Value = "27.24"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").Value = "0.0₃0850"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("D37").Value = "3.39"
$ws.Range("E37").Value = "  +4.59%  "
$ws.Range("D38").Value = "6.12"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").Value = "2.22"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").Value = "9.38"
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").Value = "445.86"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").Value = "0.285"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "0.111"
$ws.Range("E45").Value = "  +3.30%  "
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").Value = "2.803.82"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").Value = "132.27"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "25.27"
$ws.Range("E50").Value = "  +4.39%  "
$ws.Range("D51").Value = "2.27"
$ws.Range("E51").Value = "  +1.09%  "
